$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Insert a new row above row 2 (everything from the old row 2 onward shifts
# down by one), then populate the new row with the "Project" / "Google"
# attribute pair.
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "Project"
$ws.Range("B2").Value = "Google"

# The old "GmailURL" row (now row 5, column B) loses its hyperlink and the
# Hyperlink cell style/formatting that came with it.
$target = $ws.Range("B5")
$target.Hyperlinks.Delete()
$target.ClearFormats()

# The Hyperlink cell style is no longer used anywhere in the workbook now,
# so drop it from the style list too.
$wb.Styles.Item("Hyperlink").Delete()

# Leave the cursor where the author left it when they saved.
$ws.Range("F8").Select() | Out-Null
